# Update monthly rainfall rates on the "Climate" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Climate")

# rain.rate column (D) updates for Marts, April, Maj, Sommer
$ws.Range("D2").Value = 0.056
$ws.Range("D3").Value = 0.12
$ws.Range("D4").Value = 0.072
$ws.Range("D5").Value = 0.11

# Update the active cell / selection on the sheet
$ws.Activate()
$ws.Range("E8").Select()
